# Updates the cryptos list (Price / Volume(1h) columns, plus a couple of
# coin re-rankings) to match the latest scrape, as produced by the
# "Updated cryptos list ... with GitHub Actions" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.811.22'
$ws.Range("E2").Value = '  +0.67%  '

$ws.Range("D3").Value = '2.539.34'
$ws.Range("E3").Value = '  -0.26%  '

$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = '  +0.20%  '

$ws.Range("D5").Value = '311.26'
$ws.Range("E5").Value = '  +0.78%  '

$ws.Range("D6").Value = '100.77'
$ws.Range("E6").Value = '  +3.45%  '

$ws.Range("E7").Value = '  -1.01%  '

$ws.Range("E8").Value = '  +0.14%  '

$ws.Range("D9").Value = '0.523'
$ws.Range("E9").Value = '  -1.10%  '

$ws.Range("D10").Value = '35.74'
$ws.Range("E10").Value = '  +0.83%  '

$ws.Range("E11").Value = '  +0.01%  '

$ws.Range("E12").Value = '  -0.82%  '

$ws.Range("E13").Value = '  +1.27%  '

$ws.Range("D14").Value = '2.929.74'
$ws.Range("E14").Value = '  -0.19%  '

$ws.Range("D15").Value = '2.557.74'
$ws.Range("E15").Value = '  +0.31%  '

$ws.Range("E16").Value = '  -3.11%  '

$ws.Range("D17").Value = '0.816'
$ws.Range("E17").Value = '  -2.19%  '

$ws.Range("D18").Value = '42.810.62'
$ws.Range("E18").Value = '  +0.45%  '

$ws.Range("D19").Value = '6.74'
$ws.Range("E19").Value = '  +0.08%  '

$ws.Range("D20").Value = '12.34'
$ws.Range("E20").Value = '  -0.25%  '

$ws.Range("D21").Value = '0.0₃0953'
$ws.Range("E21").Value = '  -0.11%  '

$ws.Range("D22").Value = '70.17'
$ws.Range("E22").Value = '  +1.28%  '

$ws.Range("D23").Value = '243.88'
$ws.Range("E23").Value = '  -1.37%  '

$ws.Range("D24").Value = '2.88'
$ws.Range("E24").Value = '  -0.89%  '

$ws.Range("E25").Value = '  -0.57%  '

$ws.Range("E26").Value = '  -0.03%  '

$ws.Range("E27").Value = '  -4.20%  '

$ws.Range("E28").Value = '  -1.04%  '

$ws.Range("D29").Value = '10.17'
$ws.Range("E29").Value = '  +0.68%  '

$ws.Range("D30").Value = '38.69'
$ws.Range("E30").Value = '  -4.26%  '

$ws.Range("B31").Value = 'Monero'
$ws.Range("C31").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D31").Value = '158.66'
$ws.Range("E31").Value = '  +0.46%  '

$ws.Range("B32").Value = 'Filecoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D32").Value = '5.86'
$ws.Range("E32").Value = '  +2.43%  '

$ws.Range("D33").Value = '2.76'
$ws.Range("E33").Value = '  +6.93%  '

$ws.Range("E34").Value = '  +2.32%  '

$ws.Range("E35").Value = '  -0.03%  '

$ws.Range("B36").Value = 'LidoDAOToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D36").Value = '3.16'
$ws.Range("E36").Value = '  -3.33%  '

$ws.Range("B37").Value = 'Celestia'
$ws.Range("C37").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D37").Value = '18.15'
$ws.Range("E37").Value = '  -0.86%  '

$ws.Range("E38").Value = '  -5.03%  '

$ws.Range("E39").Value = '  +0.42%  '

$ws.Range("E40").Value = '  -0.13%  '

$ws.Range("D41").Value = '4.14'
$ws.Range("E41").Value = '  +2.16%  '

$ws.Range("D42").Value = '21.84'
$ws.Range("E42").Value = '  -2.62%  '

$ws.Range("E43").Value = '  +0.48%  '

$ws.Range("D44").Value = "'3.30"
$ws.Range("E44").Value = '  +3.49%  '

$ws.Range("E45").Value = '  +0.16%  '

$ws.Range("D46").Value = '1.996.97'
$ws.Range("E46").Value = '  +0.27%  '

$ws.Range("D47").Value = '9.21'
$ws.Range("E47").Value = '  +2.10%  '

$ws.Range("D48").Value = '2.783.29'
$ws.Range("E48").Value = '  -0.22%  '

$ws.Range("E49").Value = '  +0.59%  '

$ws.Range("D50").Value = '80.16'
$ws.Range("E50").Value = '  -0.87%  '

$ws.Range("D51").Value = '72.42'
$ws.Range("E51").Value = '  -1.18%  '
